$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1124, shifting existing row 1124 (and below) down to 1125.
$ws.Rows.Item(1124).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 1124 with the new record's data.
$ws.Cells.Item(1124, 1).Value  = 9
$ws.Cells.Item(1124, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1124, 3).Value  = "Metropolitana"
$ws.Cells.Item(1124, 4).Value  = 45194
$ws.Cells.Item(1124, 5).Value  = 13
$ws.Cells.Item(1124, 6).Value  = "Fruta"
$ws.Cells.Item(1124, 7).Value  = 100104
$ws.Cells.Item(1124, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(1124, 9).Value  = 100104005
$ws.Cells.Item(1124, 10).Value = "Pera"
$ws.Cells.Item(1124, 11).Value = "Packham's Triumph"
$ws.Cells.Item(1124, 12).Value = "Calibre 80"
$ws.Cells.Item(1124, 13).Value = 490
$ws.Cells.Item(1124, 14).Value = 17000
$ws.Cells.Item(1124, 15).Value = 17000
$ws.Cells.Item(1124, 16).Value = 17000
$ws.Cells.Item(1124, 17).Value = "$/caja 18 kilos embalada"
$ws.Cells.Item(1124, 18).Value = "Provincia de Linares"
$ws.Cells.Item(1124, 19).Value = 944
$ws.Cells.Item(1124, 20).Value = 18
